$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 135, pushing the existing rows 135-141
# (the previous week's price data) down to 139-145.
$ws.Rows("135:138").Insert()

# New weekly price data (week of 2023-10-24, serial 45223) for
# Agricola del Norte S.A. de Arica - Frutilla, in the usual
# Especial/Primera/Segunda/Tercera quality order.
$newRows = @(
    @{ Row = 135; Calidad = "Especial"; Volumen = 180; Min = 7000; Max = 8000; Prom = 7500; PrecioKg = 2500 },
    @{ Row = 136; Calidad = "Primera";  Volumen = 250; Min = 5000; Max = 6000; Prom = 5600; PrecioKg = 1867 },
    @{ Row = 137; Calidad = "Segunda";  Volumen = 200; Min = 4000; Max = 5000; Prom = 4500; PrecioKg = 1500 },
    @{ Row = 138; Calidad = "Tercera";  Volumen = 130; Min = 3000; Max = 4000; Prom = 3615; PrecioKg = 1205 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($row, 4).Value = 45223
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = "$/bandeja 3 kilos"
    $ws.Cells.Item($row, 18).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = 3
}
